$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2-173 hold a "last updated" date that moved forward
# by one day (2023-10-06 -> 2023-10-07, i.e. serial 45205 -> 45206).
$range = $ws.Range("C2:C173")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
